$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the contents of row 5 and row 6 for the columns that
# differ between the two records (A, B, D, E, F, G, H, S, AC, AI).
# Columns C, P, Q, R, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY
# already hold identical values in both rows, so they are left untouched.

# --- Capture current ("before") values for row 5 ---
$A5 = $ws.Range("A5").Value2
$B5 = $ws.Range("B5").Value2
$D5 = $ws.Range("D5").Value2
$E5 = $ws.Range("E5").Value2
$F5 = $ws.Range("F5").Value2
$G5 = $ws.Range("G5").Value2
$H5 = $ws.Range("H5").Value2
$S5 = $ws.Range("S5").Value2
$AC5 = $ws.Range("AC5").Value2
$AI5 = $ws.Range("AI5").Value2

# --- Capture current ("before") values for row 6 ---
$A6 = $ws.Range("A6").Value2
$B6 = $ws.Range("B6").Value2
$D6 = $ws.Range("D6").Value2
$E6 = $ws.Range("E6").Value2
$F6 = $ws.Range("F6").Value2
$G6 = $ws.Range("G6").Value2
$H6 = $ws.Range("H6").Value2
$S6 = $ws.Range("S6").Value2
$AC6 = $ws.Range("AC6").Value2
$AI6 = $ws.Range("AI6").Value2

# --- Write row 6's former values into row 5 ---
$ws.Range("A5").Value2 = $A6
$ws.Range("B5").Value2 = $B6
$ws.Range("D5").Value2 = $D6
$ws.Range("E5").Value2 = $E6
$ws.Range("F5").Value2 = $F6
$ws.Range("G5").Value2 = $G6
$ws.Range("H5").Value2 = $H6
$ws.Range("S5").Value2 = $S6
$ws.Range("AC5").Value2 = $AC6
$ws.Range("AI5").Value2 = $AI6

# --- Write row 5's former values into row 6 ---
$ws.Range("A6").Value2 = $A5
$ws.Range("B6").Value2 = $B5
$ws.Range("D6").Value2 = $D5
$ws.Range("E6").Value2 = $E5
$ws.Range("F6").Value2 = $F5
$ws.Range("G6").Value2 = $G5
$ws.Range("H6").Value2 = $H5
$ws.Range("S6").Value2 = $S5
$ws.Range("AC6").Value2 = $null
$ws.Range("AI6").Value2 = $null
